# Updates cryptos list price (D) and volume/1h (E) columns per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '43.765.51' },
    @{ Cell = 'E2'; Value = '  +0.55%  ' },
    @{ Cell = 'D3'; Value = '2.317.26' },
    @{ Cell = 'E3'; Value = '  +4.23%  ' },
    @{ Cell = 'E4'; Value = '  -0.04%  ' },
    @{ Cell = 'D5'; Value = '''97.78' },
    @{ Cell = 'E5'; Value = '  +6.12%  ' },
    @{ Cell = 'D6'; Value = '''270.94' },
    @{ Cell = 'E6'; Value = '  +0.32%  ' },
    @{ Cell = 'D7'; Value = '''0.625' },
    @{ Cell = 'E7'; Value = '  +0.52%  ' },
    @{ Cell = 'E8'; Value = '  -0.05%  ' },
    @{ Cell = 'E9'; Value = '  +0.75%  ' },
    @{ Cell = 'D10'; Value = '''45.73' },
    @{ Cell = 'E10'; Value = '  +0.19%  ' },
    @{ Cell = 'D11'; Value = '''0.0947' },
    @{ Cell = 'E11'; Value = '  +2.21%  ' },
    @{ Cell = 'D12'; Value = '''8.13' },
    @{ Cell = 'E12'; Value = '  -1.03%  ' },
    @{ Cell = 'E13'; Value = '  +0.59%  ' },
    @{ Cell = 'D14'; Value = '2.656.23' },
    @{ Cell = 'E14'; Value = '  +3.68%  ' },
    @{ Cell = 'D15'; Value = '''15.46' },
    @{ Cell = 'E15'; Value = '  +2.88%  ' },
    @{ Cell = 'E16'; Value = '  +8.52%  ' },
    @{ Cell = 'D17'; Value = '2.326.00' },
    @{ Cell = 'E17'; Value = '  +4.32%  ' },
    @{ Cell = 'D18'; Value = '43.713.21' },
    @{ Cell = 'E18'; Value = '  +0.51%  ' },
    @{ Cell = 'E19'; Value = '  +5.43%  ' },
    @{ Cell = 'D20'; Value = '''6.41' },
    @{ Cell = 'E20'; Value = '  +7.15%  ' },
    @{ Cell = 'D21'; Value = '''72.68' },
    @{ Cell = 'E21'; Value = '  +3.40%  ' },
    @{ Cell = 'D22'; Value = '''239.64' },
    @{ Cell = 'E22'; Value = '  +3.13%  ' },
    @{ Cell = 'E23'; Value = '  -2.64%  ' },
    @{ Cell = 'D24'; Value = '''9.40' },
    @{ Cell = 'E24'; Value = '  +4.34%  ' },
    @{ Cell = 'E25'; Value = '  -0.07%  ' },
    @{ Cell = 'E26'; Value = '  +0.72%  ' },
    @{ Cell = 'D27'; Value = '''11.30' },
    @{ Cell = 'E27'; Value = '  +0.33%  ' },
    @{ Cell = 'D28'; Value = '''3.47' },
    @{ Cell = 'E28'; Value = '  -2.05%  ' },
    @{ Cell = 'E29'; Value = '  +0.79%  ' },
    @{ Cell = 'D30'; Value = '''38.07' },
    @{ Cell = 'E30'; Value = '  -6.04%  ' },
    @{ Cell = 'D31'; Value = '''22.41' },
    @{ Cell = 'E31'; Value = '  +7.85%  ' },
    @{ Cell = 'D32'; Value = '''175.19' },
    @{ Cell = 'E32'; Value = '  +1.62%  ' },
    @{ Cell = 'D33'; Value = '''0.0900' },
    @{ Cell = 'E33'; Value = '  -2.04%  ' },
    @{ Cell = 'D34'; Value = '''5.48' },
    @{ Cell = 'E34'; Value = '  +0.66%  ' },
    @{ Cell = 'E35'; Value = '  +3.34%  ' },
    @{ Cell = 'D36'; Value = '''0.0360' },
    @{ Cell = 'E36'; Value = '  +2.59%  ' },
    @{ Cell = 'E37'; Value = '  -2.82%  ' },
    @{ Cell = 'D38'; Value = '''4.39' },
    @{ Cell = 'E38'; Value = '  +2.30%  ' },
    @{ Cell = 'D39'; Value = '''3.36' },
    @{ Cell = 'E39'; Value = '  -5.17%  ' },
    @{ Cell = 'D40'; Value = '''0.244' },
    @{ Cell = 'E40'; Value = '  +12.12%  ' },
    @{ Cell = 'D41'; Value = '''2.34' },
    @{ Cell = 'E41'; Value = '  +8.63%  ' },
    @{ Cell = 'E42'; Value = '  +19.33%  ' },
    @{ Cell = 'D43'; Value = '''12.16' },
    @{ Cell = 'E43'; Value = '  -2.24%  ' },
    @{ Cell = 'D44'; Value = '''9.17' },
    @{ Cell = 'E44'; Value = '  +9.74%  ' },
    @{ Cell = 'D45'; Value = '''61.91' },
    @{ Cell = 'E45'; Value = '  -1.93%  ' },
    @{ Cell = 'D46'; Value = '''5.36' },
    @{ Cell = 'E46'; Value = '  +0.91%  ' },
    @{ Cell = 'E47'; Value = '  +4.51%  ' },
    @{ Cell = 'D48'; Value = '''100.22' },
    @{ Cell = 'E48'; Value = '  +0.16%  ' },
    @{ Cell = 'E49'; Value = '  +1.16%  ' },
    @{ Cell = 'D50'; Value = '''0.190' },
    @{ Cell = 'E50'; Value = '  +17.68%  ' },
    @{ Cell = 'D51'; Value = '2.542.93' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
